$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.316.75'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.664.59'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.17'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.663.87'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.56%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.04'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.151.17'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.50%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '72.228.52'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.24'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.662.15'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.89'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.95'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '369.73'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.28%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.04'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +9.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.29'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.33'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.88'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '517.94'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.07'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.85%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.81'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.42'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.12'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.38'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.110'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.59%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.01'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.333'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.23'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '152.31'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0768'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.57%  '
